# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-39, replacing the previously
# written Strike# based values with the recalculated K values.
$newK = @{
    2  = 0
    3  = 5
    4  = 7
    5  = 9
    6  = 5
    7  = 4
    8  = 3
    9  = 8
    10 = 2
    11 = 3
    12 = 4
    13 = 6
    14 = 7
    15 = 4
    16 = 8
    17 = 6
    18 = 4
    19 = 8
    20 = 6
    21 = 8
    22 = 6
    23 = 3
    24 = 7
    25 = 11
    26 = 7
    27 = 6
    28 = 8
    29 = 5
    30 = 6
    31 = 5
    32 = 7
    33 = 11
    34 = 7
    35 = 8
    36 = 9
    37 = 5
    38 = 3
    39 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
